$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New log rows appended after the existing last row (197).
# Columns: A=run_id, B=rss_url_id, C=date, D=response, E=item_count
$newRows = @(
    @(197, 1, "2024-06-19 04:15:40", 200, 8),
    @(198, 2, "2024-06-19 04:15:40", 200, 0)
)

$startRow = 198
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
